$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Curriculum")

# Row 9 (Week 2 topic): shorten "Agile, Testing, Debugging" -> "Agile, Testing"
$ws.Range("B9").Value = "Agile, Testing"

# Row 13 (Week 2, Day Th): "Gemfile, Git" -> "Gemfile, TDD"; add new assignment note
$ws.Range("D13").Value = "Gemfile, TDD"
$ws.Range("H13").Value = "Employees and Depts?"

# Week 3 (rows 16-19) gets broken out with new day topics
$ws.Range("D16").Value = "Databases, Migrations"
$ws.Range("F16").Value = "MVP"
$ws.Range("H16").Value = "Data Diagrams (many)"

$ws.Range("D17").Value = "ActiveRecord, Dev vs. Test"
$ws.Range("H17").Value = "Emp & Dept in Database"

$ws.Range("D18").Value = "Associations, Validations"

# Week 4 (rows 22-25) gets new day topics
$ws.Range("D22").Value = "Rails, Web, HTML Verbs"
$ws.Range("D23").Value = "Router, Controllers"
$ws.Range("D24").Value = "Regex, ActiveModel Serializers"
$ws.Range("D25").Value = "Dev vs. Prod, Heroku"
$ws.Range("H25").Value = "Make-your-own API"

# Week 5 (rows 29-31) gets new day topics
$ws.Range("D29").Value = "HTML Forms, ERB"
$ws.Range("H29").Value = "Motivational Quotations"

$ws.Range("D30").Value = "Scaffold, Helpers/Partials"
$ws.Range("H30").Value = "Wallet"

$ws.Range("H31").Value = "Health Tracker"

# Update the view's active cell in the frozen bottom-right pane
$ws.Range("F25").Select()
